# Replace Thread.Sleep() method to Explicit wait in IPAIAM module.
#
# The underlying spreadsheet change: the "Runmode" column (column D) for the
# Customercare010 .. Customercare014 test cases was flipped from "N" (skip)
# to "Y" (run), now that the tests no longer rely on a fixed Thread.Sleep()
# and use an explicit wait instead, so they can be safely re-enabled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-13 hold Customercare010 .. Customercare014; column D is "Runmode".
$ws.Range("D9").Value2  = "Y"
$ws.Range("D10").Value2 = "Y"
$ws.Range("D11").Value2 = "Y"
$ws.Range("D12").Value2 = "Y"
$ws.Range("D13").Value2 = "Y"

# Leave the sheet's on-screen selection where the author left it after
# making the edit.
$ws.Range("D8:D19").Select()
